# "Fixed typos" - add the missing "Language translation" bullet to the
# "We will apply NLG for" slide (slide 4), in the body placeholder that
# already lists "Building a chatbot" / " Essay writing".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(8)
$tr = $shp.TextFrame.TextRange

# Append a new paragraph (carriage return + text) after the existing text
# so PowerPoint creates a new <a:p> that inherits the paragraph/run
# formatting of the last paragraph, matching the other bullet items.
$newRange = $tr.InsertAfter("`rLanguage translation")
